$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-price data refresh to the per-job
# "Profits" sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). For each changed leve row
# it rewrites the price/profit columns (H:N) with freshly pulled values. A few
# rows had their optional HQ columns (M/N) added or removed because the item no
# longer has (or newly has) an HQ market price.

$ws = $wb.Worksheets.Item("ALC")
$ALCValues = @{
    "K18" = 880.1111
    "M18" = -596.1111
    "N18" = -1346.1429
    "H18" = 835.5
    "I18" = 880.1111
    "J18" = 778.1429000000001
    "L18" = 778.1429000000001
    "N70" = -7629.299999999999
    "H70" = 2281.75
    "J70" = 2363.1
    "L70" = 7089.299999999999
    "N73" = -8961.299999999999
    "H73" = 2281.75
    "J73" = 2363.1
    "L73" = 7089.299999999999
    "K138" = 7752.999899999999
    "I138" = 2584.3333
    "M138" = -2612.999899999999
    "H138" = 3618.9644
}
foreach ($cellRef in $ALCValues.Keys) {
    $ws.Range($cellRef).Value = $ALCValues[$cellRef]
}

$ws = $wb.Worksheets.Item("ARM")
$ARMValues = @{
    "J5" = 1000
    "L5" = 1000
    "K5" = 189.72728
    "I5" = 189.72728
    "M5" = -77.72728000000001
    "H5" = 257.25
    "N5" = -1224
    "K31" = 10372.25
    "I31" = 10372.25
    "M31" = -10078.25
    "H31" = 10372.25
    "M32" = -406670.2
    "I32" = 406957.2
    "H32" = 199070.8
    "K32" = 406957.2
    "M61" = -4414.9644
    "H61" = 2861657.2
    "I61" = 4626.9644
    "K61" = 4626.9644
    "N94" = -52132
    "H94" = 50330
    "J94" = 50330
    "L94" = 50330
    "I97" = 6555
    "M97" = -6059
    "H97" = 166671680
    "K97" = 6555
    "K102" = 651.4
    "M102" = 970.6
    "I102" = 651.4
    "H102" = 651.4
    "N108" = -42536.715
    "H108" = 34856.715
    "J108" = 34856.715
    "L108" = 34856.715
    "H136" = 2861657.2
    "K136" = 13880.8932
    "M136" = -11330.8932
    "I136" = 4626.9644
}
foreach ($cellRef in $ARMValues.Keys) {
    $ws.Range($cellRef).Value = $ARMValues[$cellRef]
}

$ws = $wb.Worksheets.Item("BSM")
$BSMValues = @{
    "H4" = 257.25
    "I4" = 189.72728
    "J4" = 1000
    "L4" = 1000
    "K4" = 189.72728
    "M4" = -74.72728000000001
    "N4" = -1230
    "K105" = 6411.0557
    "M105" = -4664.0557
    "I105" = 6411.0557
    "H105" = 6827.276
}
foreach ($cellRef in $BSMValues.Keys) {
    $ws.Range($cellRef).Value = $BSMValues[$cellRef]
}

$ws = $wb.Worksheets.Item("CRP")
$CRPValues = @{
    "N16" = -169388.5
    "H16" = 168552.17
    "J16" = 168814.5
    "L16" = 168814.5
    "N113" = -173154.5
    "H113" = 168552.17
    "J113" = 168814.5
    "L113" = 168814.5
    "M132" = -1717.142599999999
    "N132" = -11264
    "H132" = 1652.909
    "I132" = 1415.7142
    "J132" = 2068
    "L132" = 6204
    "K132" = 4247.142599999999
    "N141" = -240932.17
    "H141" = 211918.92
    "J141" = 230572.17
    "L141" = 230572.17
}
foreach ($cellRef in $CRPValues.Keys) {
    $ws.Range($cellRef).Value = $CRPValues[$cellRef]
}

$ws = $wb.Worksheets.Item("CUL")
$CULValues = @{
    "J31" = 0
    "L31" = 0
    "H31" = 350
    "N38" = -921.1176399999999
    "H38" = 52.52941
    "J38" = 75.70587999999999
    "L38" = 227.11764
    "J40" = 100
    "L40" = 400
    "N40" = -538
    "H40" = 43.333332
    "H62" = 2733.111
    "K62" = 8298.999899999999
    "I62" = 2766.3333
    "M62" = -7612.999899999999
    "I65" = 2766.3333
    "H65" = 2733.111
    "K65" = 24896.9997
    "M65" = -21464.9997
    "K69" = 9512.000100000001
    "M69" = -8701.000100000001
    "I69" = 3170.6667
    "H69" = 7716
    "H72" = 7716
    "K72" = 28536.0003
    "M72" = -24480.0003
    "I72" = 3170.6667
    "K82" = 15039
    "M82" = -14633
    "H82" = 15001.733
    "I82" = 5013
    "K85" = 15039
    "M85" = -13635
    "H85" = 15001.733
    "I85" = 5013
    "M87" = -12354.8568
    "H87" = 11811.667
    "I87" = 4534.2856
    "K87" = 13602.8568
    "I90" = 4534.2856
    "M90" = -34568.5704
    "H90" = 11811.667
    "K90" = 40808.5704
    "L92" = 876
    "N92" = -3372
    "H92" = 274.875
    "J92" = 292
    "J98" = 0
    "L98" = 0
    "K98" = 0
    "I98" = 0
    "H98" = 0
    "I104" = 2633.3333
    "J104" = 3998.6667
    "L104" = 11996.0001
    "K104" = 7899.999899999999
    "N104" = -17238.0001
    "M104" = -5278.999899999999
    "H104" = 3316
    "J127" = 7891.25
    "L127" = 23673.75
    "N127" = -33593.75
    "H127" = 7891.25
    "N132" = -134786
    "H132" = 8742.857
    "J132" = 14414
    "L132" = 129726
}
foreach ($cellRef in $CULValues.Keys) {
    $ws.Range($cellRef).Value = $CULValues[$cellRef]
}
foreach ($cellRef in @("N31","M98","N98")) {
    $ws.Range($cellRef).ClearContents()
}

$ws = $wb.Worksheets.Item("GSM")
$GSMValues = @{
    "K80" = 2734
    "M80" = -1736
    "H80" = 2734
    "I80" = 2734
    "K83" = 13670
    "M83" = -8678
    "I83" = 2734
    "H83" = 2734
    "N138" = -109278
    "H138" = 98998
    "J138" = 98998
    "L138" = 98998
}
foreach ($cellRef in $GSMValues.Keys) {
    $ws.Range($cellRef).Value = $GSMValues[$cellRef]
}

$ws = $wb.Worksheets.Item("LTW")
$LTWValues = @{
    "J40" = 9992.5
    "L40" = 9992.5
    "K40" = 5425.421
    "M40" = -5289.421
    "N40" = -10264.5
    "H40" = 5860.381
    "I40" = 5425.421
    "N94" = -21517
    "H94" = 20165
    "J94" = 20165
    "L94" = 20165
    "K99" = 30259
    "M99" = -27264
    "H99" = 30259
    "I99" = 30259
    "M122" = -6949.599999999999
    "I122" = 3133.2
    "H122" = 3312.375
    "K122" = 9399.599999999999
}
foreach ($cellRef in $LTWValues.Keys) {
    $ws.Range($cellRef).Value = $LTWValues[$cellRef]
}

$ws = $wb.Worksheets.Item("WVR")
$WVRValues = @{
    "N113" = -10928.9999
    "H113" = 2169.3157
    "I113" = 2156.8462
    "J113" = 2196.3333
    "L113" = 6588.999899999999
    "K113" = 6470.5386
    "M113" = -4300.5386
    "M122" = -13471.4995
    "N122" = -851661.25
    "I122" = 5307.1665
    "H122" = 74543.81
    "J122" = 282253.75
    "L122" = 846761.25
    "K122" = 15921.4995
}
foreach ($cellRef in $WVRValues.Keys) {
    $ws.Range($cellRef).Value = $WVRValues[$cellRef]
}
